$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two outlier rows (original rows 2 and 4, identified by the
# duplicated/erroneous weight value 10.90219764355227) and shift the
# remaining rows up, producing the final evaluation dataset.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(2).Delete()
